$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds numeric-looking text (e.g. "1.0000", "28.818.70").
# Force the affected cells to Text format before writing so Excel doesn't
# silently coerce the string into a Number and drop significant trailing
# zeros / thousands-style dot groupings.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.818.70'
$ws.Range("E2").Value = '  +7.91%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.809.63'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.0000'
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.64'
$ws.Range("E5").Value = '  +3.26%  '
$ws.Range("E6").Value = '  +0.15%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4969'
$ws.Range("E7").Value = '  +3.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '43.18'
$ws.Range("E8").Value = '  +3.97%  '
$ws.Range("E9").Value = '  +7.99%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06414'
$ws.Range("E10").Value = '  +3.74%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.814.64'
$ws.Range("E11").Value = '  +5.38%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '16.76'
$ws.Range("E12").Value = '  +5.66%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07073'
$ws.Range("E13").Value = '  +3.37%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6457'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '84.04'
$ws.Range("E15").Value = '  +9.43%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.692'
$ws.Range("E16").Value = '  +5.23%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '28.810.67'
$ws.Range("E17").Value = '  +8.56%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9999'
$ws.Range("E18").Value = '  +0.12%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007334'
$ws.Range("E19").Value = '  +2.65%  '
$ws.Range("E20").Value = '  +0.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.24'
$ws.Range("E21").Value = '  +7.93%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.043.71'
$ws.Range("E22").Value = '  +4.99%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.578'
$ws.Range("E23").Value = '  +3.71%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.882'
$ws.Range("E24").Value = '  +3.65%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.336'
$ws.Range("E25").Value = '  +5.63%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '142.34'
$ws.Range("E26").Value = '  +2.24%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '128.63'
$ws.Range("E27").Value = '  +21.37%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.36'
$ws.Range("E28").Value = '  +7.57%  '
$ws.Range("E29").Value = '  +6.48%  '
$ws.Range("E30").Value = '  +3.43%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.127'
$ws.Range("E31").Value = '  +2.63%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08344'
$ws.Range("E32").Value = '  +5.64%  '
$ws.Range("E33").Value = '  +3.76%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04945'
$ws.Range("E34").Value = '  +10.45%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.090'
$ws.Range("E35").Value = '  +9.18%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.710'
$ws.Range("E36").Value = '  +4.41%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6690'
$ws.Range("E37").Value = '  +8.58%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.272'
$ws.Range("E38").Value = '  +13.35%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.729'
$ws.Range("E39").Value = '  +12.00%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9542'
$ws.Range("E40").Value = '  +3.17%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.088'
$ws.Range("E41").Value = '  +8.78%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.01590'
$ws.Range("E42").Value = '  +6.88%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9999'
$ws.Range("E43").Value = '  +0.20%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4078'
$ws.Range("E44").Value = '  +6.72%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '99.77'
$ws.Range("E45").Value = '  -0.12%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.143'
$ws.Range("E46").Value = '  +5.51%  '
$ws.Range("E47").Value = '  +5.95%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05510'
$ws.Range("E48").Value = '  +2.65%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '31.57'
$ws.Range("E49").Value = '  +5.23%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.057'
$ws.Range("E50").Value = '  +2.80%  '
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.307'
$ws.Range("E51").Value = '  +5.82%  '
